$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.745.38'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.522.19'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''608.11'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '''196.12'
$ws.Range('E6').Value = '  +5.36%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  -6.64%  '
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('D11').Value = '''53.74'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('D13').Value = '''9.50'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').Value = '4.082.37'
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').Value = '''596.76'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = '''12.82'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').Value = '69.919.12'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').Value = '3.533.33'
$ws.Range('E19').Value = '  +1.72%  '
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').Value = '''18.07'
$ws.Range('E22').Value = '  +5.47%  '
$ws.Range('D23').Value = '''5.32'
$ws.Range('E23').Value = '  +5.42%  '
$ws.Range('D24').Value = '''102.32'
$ws.Range('E24').Value = '  -2.92%  '
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  +4.60%  '
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').Value = '''33.50'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('D30').Value = '''7.05'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').Value = '''4.23'
$ws.Range('E31').Value = '  +3.85%  '
$ws.Range('D32').Value = '''12.44'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '''63.16'
$ws.Range('E34').Value = '  -0.14%  '
$ws.Range('D35').Value = '0.0₃0850'
$ws.Range('E35').Value = '  +9.09%  '
$ws.Range('D36').Value = '3.713.38'
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').Value = '''0.392'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('D41').Value = '''36.55'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('D42').Value = '''491.15'
$ws.Range('E42').Value = '  -6.42%  '
$ws.Range('D43').Value = '''0.132'
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('D44').Value = '''0.0454'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('E46').Value = '  -4.36%  '
$ws.Range('D47').Value = '''3.28'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('D49').Value = '''8.49'
$ws.Range('E49').Value = '  -3.58%  '
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('E51').Value = '  +10.95%  '
